$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "username" / "pswd" columns (D/E) mirroring the customer name (A)
# and a constant placeholder password, with a couple of name corrections
# applied to the new username column (row 3: Sawik -> Satwik, row 19:
# Heemmanshuu -> Heemmanshu).

$ws.Range("D1").Value = "username"
$ws.Range("E1").Value = "pswd"

$usernames = @{
    2  = "Thivesh"
    3  = "Satwik"
    4  = "Akash"
    5  = "Karthikeya"
    6  = "Vishnu"
    7  = "Aadarsh"
    8  = "Abhijith"
    9  = "Aditya"
    10 = "Saathwick"
    11 = "Kranthi"
    12 = "Likith"
    13 = "Ashwin"
    14 = "Bhanu"
    15 = "Sekhar"
    16 = "Bharat"
    17 = "Varma"
    18 = "Gattu"
    19 = "Heemmanshu"
    20 = "Hitesh"
    21 = "Jayanth"
}

foreach ($row in 2..21) {
    $ws.Range("D$row").Value = $usernames[$row]
    $ws.Range("E$row").Value = "abcd"
}

$ws.Range("D21").Select()
